$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Delete the "License Information" Heading2 paragraph (originally
#    paragraph #4).
# ------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.Delete()

# ------------------------------------------------------------------
# 2. Rewrite the paragraph that starts with the bold
#    "मुख्य शब्द (Biblica)" run (now paragraph #4 after the delete
#    above).  Keep the leading empty run, replace everything else
#    with the new text/runs, and keep a trailing empty run.
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(4)
$rng = $p.Range
$rng.MoveEnd(1, -1)      # exclude the paragraph mark
$rng.Text = ""           # clear paragraph content (keeps it a Normal para)
$start = $rng.Start

$seg1 = "Biblica Study Notes (Key Terms)"
$seg2 = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
$seg3 = "Biblica Study Notes"
$seg4 = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."

$rng.InsertAfter($seg1)
$rng.InsertAfter($seg2)
$rng.InsertAfter($seg3)
$rng.InsertAfter($seg4)

$p1start = $start
$p1end = $p1start + $seg1.Length
$p2start = $p1end
$p2end = $p2start + $seg2.Length
$p3start = $p2end
$p3end = $p3start + $seg3.Length
$p4start = $p3end
$p4end = $p4start + $seg4.Length

$r1 = $d.Range($p1start, $p1end)
$r1.Bold = 1

$r2 = $d.Range($p2start, $p2end)
$r2.Bold = 1
$r2.Bold = 0

$r3 = $d.Range($p3start, $p3end)
$r3.Bold = 1
$r3.Bold = 0

$r4 = $d.Range($p4start, $p4end)
$r4.Bold = 1
$r4.Bold = 0

# ------------------------------------------------------------------
# 3. Delete the "This PDF version is provided under the same
#    license." paragraph (now paragraph #5).
# ------------------------------------------------------------------
$d.Paragraphs.Item(5).Range.Delete()

# ------------------------------------------------------------------
# 4. Delete the paragraph holding the italic keyword list that
#    follows the "पर" Heading2 paragraph.
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute("प्रकाश, प्रभु का दिन", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $kwPara = $d.Content.Find.Parent.Paragraphs.Item(1)
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.StartsWith("प्रकाश, प्रभु का दिन")) {
        $pp.Range.Delete()
        break
    }
}
